$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Model Accuracy (-0.3, 0.3, 0.3)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Model Accuracy (-0.3, 0.3, 0.3)")

# New header cells C1:G1 (use same style/format as B1 header cell)
$ws1.Range("C1").Value = "Market threshold"
$ws1.Range("D1").Value = "Market min"
$ws1.Range("E1").Value = "Market max"
$ws1.Range("F1").Value = "Recall"
$ws1.Range("G1").Value = "Precision"
$ws1.Range("B1").Copy()
$ws1.Range("C1:G1").PasteSpecial(-4122)

# Updated accuracy column + new metric columns
$ws1.Range("B2").Value = 56.47921760391198
$ws1.Range("C2").Value = 0.05450546436368681
$ws1.Range("D2").Value = -15.55441
$ws1.Range("E2").Value = 15.06418
$ws1.Range("F2").Value = 11.11111111111111
$ws1.Range("G2").Value = 2.083333333333333

$ws1.Range("B3").Value = 32.3960880195599
$ws1.Range("C3").Value = 0.009583939973006913
$ws1.Range("D3").Value = -19.35264
$ws1.Range("E3").Value = 13.70093
$ws1.Range("F3").Value = 8.042895442359249
$ws1.Range("G3").Value = 23.62204724409449

$ws1.Range("B4").Value = 84.65770171149144
$ws1.Range("C4").Value = 0.04158117063764853
$ws1.Range("D4").Value = -18.75314
$ws1.Range("E4").Value = 23.33066
$ws1.Range("F4").Value = 0
$ws1.Range("G4").Value = 0

$ws1.Range("B5").Value = 71.63814180929096
$ws1.Range("C5").Value = 0.02983403801513819
$ws1.Range("D5").Value = -12.78028
$ws1.Range("E5").Value = 12.42348
$ws1.Range("F5").Value = 3.636363636363636
$ws1.Range("G5").Value = 6.896551724137931

$ws1.Range("B6").Value = 85.08557457212714
$ws1.Range("C6").Value = 0.08368817696170747
$ws1.Range("D6").Value = -16.47904
$ws1.Range("E6").Value = 14.94325
$ws1.Range("F6").Value = 0
$ws1.Range("G6").Value = 0

# ---------------------------------------------------------------------------
# Sheet 2: Confusion Matrix TOTALENERGIES SE (-0.3, 0.3, 0.3)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Confusion Matrix TOTALENERGIES SE (-0.3, 0.3, 0.3)")
$ws2.Range("B3").Value = 7
$ws2.Range("C3").Value = 919
$ws2.Range("D3").Value = 6

# ---------------------------------------------------------------------------
# Sheet 3: Confusion Matrix FMC CORP (-0.3, 0.3, 0.3)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Confusion Matrix FMC CORP (-0.3, 0.3, 0.3)")
$ws3.Range("B2").Value = 30
$ws3.Range("C2").Value = 67
$ws3.Range("D2").Value = 30

$ws3.Range("B3").Value = 212
$ws3.Range("C3").Value = 379
$ws3.Range("D3").Value = 205

$ws3.Range("B4").Value = 131
$ws3.Range("C4").Value = 206
$ws3.Range("D4").Value = 121

# ---------------------------------------------------------------------------
# Sheet 4: Confusion Matrix BP PLC (-0.3, 0.3, 0.3)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Confusion Matrix BP PLC (-0.3, 0.3, 0.3)")
$ws4.Range("B3").Value = 36
$ws4.Range("C3").Value = 1377
$ws4.Range("D3").Value = 34

$ws4.Range("B4").Value = 4
$ws4.Range("C4").Value = 138
$ws4.Range("D4").Value = 8

# ---------------------------------------------------------------------------
# Sheet 5: Confusion Matrix STORA ENSO (-0.3, 0.3, 0.3)
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Confusion Matrix STORA ENSO (-0.3, 0.3, 0.3)")
$ws5.Range("B2").Value = 4
$ws5.Range("C2").Value = 53

$ws5.Range("B3").Value = 91
$ws5.Range("C3").Value = 1161
$ws5.Range("D3").Value = 99

$ws5.Range("B4").Value = 15
$ws5.Range("C4").Value = 144
$ws5.Range("D4").Value = 7

# ---------------------------------------------------------------------------
# Sheet 6: Confusion Matrix BHP GROUP (-0.3, 0.3, 0.3)
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("Confusion Matrix BHP GROUP (-0.3, 0.3, 0.3)")
$ws6.Range("B2").Value = 0
$ws6.Range("C2").Value = 116

$ws6.Range("B3").Value = 4
$ws6.Range("C3").Value = 1392
$ws6.Range("D3").Value = 2

$ws6.Range("B4").Value = 0
$ws6.Range("C4").Value = 65
